$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.266.44'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.623.95'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.51%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.36'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.108'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.377'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.25%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.52'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.58%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.152'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.23'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.090.78'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.74%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.114.85'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000146'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.89%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.643.74'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.53'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.52'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.14%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '342.60'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.17%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.88'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('B23').Value = 'LEO'
$ws.Range('C23').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.72'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.91%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.09'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.24%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.69'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.53%  '
$ws.Range('B26').Value = 'SuiNetwork'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.62'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.07%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.03'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.24%  '
$ws.Range('B28').Value = 'Bittensor'
$ws.Range('C28').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '546.98'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.11%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.162'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.42%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.99'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.43%  '
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.03'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.31%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.76'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.14%  '
$ws.Range('B34').Value = 'PEPE'
$ws.Range('C34').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0838'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.64%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.24'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.47%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '169.14'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.54%  '
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.402'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.74%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.92'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.58%  '
$ws.Range('B40').Value = 'EthereumClassic'
$ws.Range('C40').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.96'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.24%  '
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '165.27'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.06%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.90'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.75'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.40%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.86'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.98%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0563'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.624'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.39%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0243'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.67%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.94'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +12.84%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0952'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.89%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.57'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.06%  '
